$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (the rows referencing "Resolving-Mac")
$ws.Rows("4:5").Delete()

# Update row 2 with new values
$ws.Range("B2").Value = "Tac1"
$ws.Range("C2").Value = "Tacr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.190825
$ws.Range("N2").Value = 0.572475
$ws.Range("O2").Value = 0.6793270274792366
$ws.Range("P2").Value = 0.6793270274792366
$ws.Range("Q2").Value = 0.9809923966999999
$ws.Range("R2").Value = 8.8289315703
$ws.Range("S2").Value = 0.6793270274792366
$ws.Range("T2").Value = 0.6793270274792366

# Update row 3 with new values
$ws.Range("B3").Value = "Tac1"
$ws.Range("C3").Value = "Tacr2"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 0.09007799999999999
$ws.Range("O3").Value = 0.3206729725207633
$ws.Range("P3").Value = 0.3206729725207634
$ws.Range("Q3").Value = 0.4630726220879999
$ws.Range("S3").Value = 0.3206729725207633
$ws.Range("T3").Value = 0.3206729725207634
